$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text.Replace('"', '""') + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

Set-TextValue $ws.Range("D2") '62.285.65'
Set-TextValue $ws.Range("E2") '  -1.79%  '
Set-TextValue $ws.Range("D3") '2.443.43'
Set-TextValue $ws.Range("E3") '  -0.27%  '
Set-TextValue $ws.Range("E4") '  -0.05%  '
Set-TextValue $ws.Range("D5") '585.28'
Set-TextValue $ws.Range("E5") '  +2.24%  '
Set-TextValue $ws.Range("D6") '143.66'
Set-TextValue $ws.Range("E6") '  -2.22%  '
Set-TextValue $ws.Range("E7") '  +0.03%  '
Set-TextValue $ws.Range("E8") '  -1.89%  '
Set-TextValue $ws.Range("D9") '2.441.26'
Set-TextValue $ws.Range("E9") '  -0.24%  '
Set-TextValue $ws.Range("E10") '  -3.33%  '
Set-TextValue $ws.Range("E11") '  +2.18%  '
Set-TextValue $ws.Range("D12") '5.20'
Set-TextValue $ws.Range("E12") '  -0.87%  '
Set-TextValue $ws.Range("D13") '0.344'
Set-TextValue $ws.Range("E13") '  -3.48%  '
Set-TextValue $ws.Range("D14") '26.43'
Set-TextValue $ws.Range("E14") '  -2.56%  '
Set-TextValue $ws.Range("E15") '  -3.96%  '
Set-TextValue $ws.Range("D16") '2.872.34'
Set-TextValue $ws.Range("E16") '  -0.14%  '
Set-TextValue $ws.Range("D17") '62.150.59'
Set-TextValue $ws.Range("E17") '  -1.41%  '
Set-TextValue $ws.Range("D18") '2.439.97'
Set-TextValue $ws.Range("E18") '  -0.05%  '
Set-TextValue $ws.Range("D19") '10.89'
Set-TextValue $ws.Range("E19") '  -3.88%  '
Set-TextValue $ws.Range("D20") '7.11'
Set-TextValue $ws.Range("E20") '  -2.72%  '
Set-TextValue $ws.Range("D21") '330.62'
Set-TextValue $ws.Range("E21") '  +0.68%  '
Set-TextValue $ws.Range("E22") '  -2.19%  '
Set-TextValue $ws.Range("D23") '1.98'
Set-TextValue $ws.Range("E23") '  -5.98%  '
Set-TextValue $ws.Range("E24") '  +0.13%  '
Set-TextValue $ws.Range("D25") '65.77'
Set-TextValue $ws.Range("E25") '  +0.45%  '
Set-TextValue $ws.Range("D26") '9.38'
Set-TextValue $ws.Range("E26") '  +3.91%  '
Set-TextValue $ws.Range("D27") '618.09'
Set-TextValue $ws.Range("E27") '  -0.21%  '
Set-TextValue $ws.Range("D29") '0.0₃0955'
Set-TextValue $ws.Range("E29") '  -7.75%  '
Set-TextValue $ws.Range("D30") '0.999'
Set-TextValue $ws.Range("E30") '  -0.13%  '
Set-TextValue $ws.Range("E31") '  -4.96%  '
Set-TextValue $ws.Range("D32") '8.00'
Set-TextValue $ws.Range("E32") '  -3.55%  '
Set-TextValue $ws.Range("E33") '  -0.15%  '
Set-TextValue $ws.Range("E34") '  -0.62%  '
Set-TextValue $ws.Range("D35") '4.92'
Set-TextValue $ws.Range("E35") '  -6.09%  '
Set-TextValue $ws.Range("E36") '  +0.13%  '
Set-TextValue $ws.Range("D37") '1.43'
Set-TextValue $ws.Range("E37") '  -6.37%  '
Set-TextValue $ws.Range("E38") '  -1.06%  '
Set-TextValue $ws.Range("D39") '151.44'
Set-TextValue $ws.Range("E39") '  +3.89%  '
Set-TextValue $ws.Range("E40") '  -2.30%  '
Set-TextValue $ws.Range("D41") '5.24'
Set-TextValue $ws.Range("E41") '  -3.81%  '
Set-TextValue $ws.Range("E42") '  -2.12%  '
Set-TextValue $ws.Range("D43") '42.50'
Set-TextValue $ws.Range("E43") '  +1.37%  '
Set-TextValue $ws.Range("E44") '  +0.00%  '
Set-TextValue $ws.Range("E45") '  -8.66%  '
Set-TextValue $ws.Range("D46") '143.52'
Set-TextValue $ws.Range("E46") '  -3.61%  '
Set-TextValue $ws.Range("E47") '  -3.68%  '
Set-TextValue $ws.Range("E48") '  -1.70%  '
Set-TextValue $ws.Range("D49") '0.598'
Set-TextValue $ws.Range("E49") '  -0.68%  '
Set-TextValue $ws.Range("D50") '19.47'
Set-TextValue $ws.Range("E50") '  -8.46%  '
Set-TextValue $ws.Range("E51") '  -1.22%  '

$excel.CutCopyMode = $false
